# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '22.376.55'
$ws.Range('E2').Value = '  -4.22%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.568.94'
$ws.Range('E3').Value = '  -3.68%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '289.39'
$ws.Range('E6').Value = '  -2.95%  '
$ws.Range('E7').Value = '  -1.91%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '49.32'
$ws.Range('E8').Value = '  -2.01%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3379'
$ws.Range('E9').Value = '  -3.04%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.162'
$ws.Range('E10').Value = '  -3.22%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07620'
$ws.Range('E11').Value = '  -5.05%  '
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '21.26'
$ws.Range('E13').Value = '  -2.48%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.057'
$ws.Range('E14').Value = '  -3.96%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.898'
$ws.Range('E15').Value = '  -4.43%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.571.32'
$ws.Range('E16').Value = '  -3.68%  '
$ws.Range('E17').Value = '  -5.11%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '89.61'
$ws.Range('E18').Value = '  -5.31%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06748'
$ws.Range('E19').Value = '  -2.68%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.002'
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.221'
$ws.Range('E21').Value = '  -6.01%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.5346'
$ws.Range('E22').Value = '  -5.62%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '16.50'
$ws.Range('E23').Value = '  -4.64%  '
$ws.Range('E24').Value = '  -3.13%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '22.400.37'
$ws.Range('E25').Value = '  -4.18%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.376'
$ws.Range('E26').Value = '  -1.84%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.897'
$ws.Range('E27').Value = '  -2.15%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '20.00'
$ws.Range('E28').Value = '  -3.46%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '146.04'
$ws.Range('E29').Value = '  -4.05%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.958'
$ws.Range('E30').Value = '  -4.20%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '125.47'
$ws.Range('E31').Value = '  -4.51%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.748.76'
$ws.Range('E32').Value = '  -3.51%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.039'
$ws.Range('E33').Value = '  +7.19%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.247'
$ws.Range('E34').Value = '  -7.54%  '
$ws.Range('E35').Value = '  -5.30%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '10.17'
$ws.Range('E36').Value = '  -8.38%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.08469'
$ws.Range('E37').Value = '  -2.77%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02538'
$ws.Range('E38').Value = '  -5.10%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.2323'
$ws.Range('E39').Value = '  -4.08%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.06532'
$ws.Range('E40').Value = '  -2.34%  '
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.536'
$ws.Range('E41').Value = '  -5.21%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.260'
$ws.Range('E42').Value = '  -1.66%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '11.69'
$ws.Range('E43').Value = '  -7.81%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.6355'
$ws.Range('E44').Value = '  -6.46%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '14.22'
$ws.Range('E45').Value = '  -6.78%  '
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5980'
$ws.Range('E47').Value = '  -5.13%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.748'
$ws.Range('E48').Value = '  -3.64%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.107'
$ws.Range('E49').Value = '  -5.47%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.260'
$ws.Range('E50').Value = '  +3.82%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '124.75'
$ws.Range('E51').Value = '  -1.40%  '
